# Update countries & provincias Spain
#
# paises.xlsx ("Pais" sheet) refresh:
#   - A handful of country-name rows got re-paired with their correct
#     statistics (the shared-string table had Chile/Japon,
#     Panama/Republica Dominicana, Eslovenia/Azerbaiyan and Mozambique
#     shuffled relative to the per-row case counts), so those rows need
#     both their country label (col A) and numbers (cols B:H) fixed up.
#   - Several other rows just received refreshed case counts (daily
#     COVID-19 data update) with no label change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 4: Estados Unidos -- refreshed totals -----------------------
$ws.Range("B4").Value = 680541
$ws.Range("C4").Value = 2971
$ws.Range("D4").Value = 58066
$ws.Range("E4").Value = 587752
$ws.Range("G4").Value = 106
$ws.Range("H4").Value = 34723

# --- Row 12: Turquia -- refreshed totals ------------------------------
$ws.Range("B12").Value = 78546
$ws.Range("C12").Value = 4353
$ws.Range("D12").Value = 8631
$ws.Range("E12").Value = 68146
$ws.Range("F12").Value = 1845
$ws.Range("G12").Value = 126
$ws.Range("H12").Value = 1769

# --- Rows 27-28: Japon / Chile swapped, each with its own data -------
$ws.Range("A27").Value = "Chile"
$ws.Range("B27").Value = 9252
$ws.Range("C27").Value = 445
$ws.Range("D27").Value = 3299
$ws.Range("E27").Value = 5837
$ws.Range("F27").Value = 384
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = 116

$ws.Range("A28").Value = "Japon"
$ws.Range("B28").Value = 9231
$ws.Range("D28").Value = 935
$ws.Range("E28").Value = 8106
$ws.Range("F28").Value = 193
$ws.Range("H28").Value = 190

# --- Row 36: Australia -- refreshed totals ----------------------------
$ws.Range("B36").Value = 6526
$ws.Range("C36").Value = 58
$ws.Range("D36").Value = 3821
$ws.Range("E36").Value = 2640

# --- Rows 48-49: Panama / Republica Dominicana swapped ----------------
$ws.Range("A48").Value = "Republica Dominicana"
$ws.Range("B48").Value = 4126
$ws.Range("C48").Value = 371
$ws.Range("D48").Value = 215
$ws.Range("E48").Value = 3711
$ws.Range("F48").Value = 146
$ws.Range("G48").Value = 4
$ws.Range("H48").Value = 200

$ws.Range("A49").Value = "Panama"
$ws.Range("B49").Value = 4016
$ws.Range("C49").Value = 265
$ws.Range("D49").Value = 98
$ws.Range("E49").Value = 3809
$ws.Range("F49").Value = 99
$ws.Range("G49").Value = 6
$ws.Range("H49").Value = 109

# --- Rows 72-73: Eslovenia / Azerbaiyan swapped -----------------------
$ws.Range("A72").Value = "Azerbaiyan"
$ws.Range("B72").Value = 1340
$ws.Range("C72").Value = 57
$ws.Range("D72").Value = 528
$ws.Range("E72").Value = 797
$ws.Range("F72").Value = 26
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 15

$ws.Range("A73").Value = "Eslovenia"
$ws.Range("B73").Value = 1304
$ws.Range("C73").Value = 36
$ws.Range("D73").Value = 174
$ws.Range("E73").Value = 1064
$ws.Range("G73").Value = 5
$ws.Range("H73").Value = 66

# --- Row 112: Mauricio -- refreshed totals ----------------------------
$ws.Range("D112").Value = 108
$ws.Range("E112").Value = 207

# --- Row 127: Tanzania -- refreshed totals ----------------------------
$ws.Range("E127").Value = 131
$ws.Range("F127").Value = 4
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 5

# --- Rows 166-169: Mozambique moved ahead of Siria/Guam/Sudan ---------
$ws.Range("A166").Value = "Mozambique"
$ws.Range("B166").Value = 34
$ws.Range("C166").Value = 3
$ws.Range("D166").Value = 2
$ws.Range("E166").Value = 32
$ws.Range("H166").Value = 0

$ws.Range("A167").Value = "Siria"
$ws.Range("B167").Value = 33
$ws.Range("D167").Value = 5
$ws.Range("E167").Value = 26
$ws.Range("H167").Value = 2

$ws.Range("A168").Value = "Guam"
$ws.Range("D168").Value = 0
$ws.Range("E168").Value = 31
$ws.Range("H168").Value = 1

$ws.Range("A169").Value = "Sudan"
$ws.Range("B169").Value = 32
$ws.Range("D169").Value = 4
$ws.Range("E169").Value = 23
$ws.Range("H169").Value = 5
